$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Clear old content fully before rewriting layout
$ws.Range("A1:F18").ClearContents()

# Column widths
$ws.Columns.Item(1).ColumnWidth = 25.83203125
$ws.Columns.Item(2).ColumnWidth = 27.1640625
$ws.Columns.Item(3).ColumnWidth = 50
$ws.Columns.Item(4).ColumnWidth = 31.1640625
$ws.Columns.Item(5).ColumnWidth = 28
$ws.Columns.Item(6).ColumnWidth = 10.83203125

# Wrap text for columns B through F (applies the new style used across the data area)
$ws.Range("B1:F18").WrapText = $true

# Cell values
$ws.Range("A1").Value = 'Category'
$ws.Range("B1").Value = 'Data '
$ws.Range("C1").Value = 'Use'
$ws.Range("D1").Value = 'Timestep'
$ws.Range("E1").Value = 'Source'
$ws.Range("A2").Value = 'Demand'
$ws.Range("B2").Value = 'Temperature'
$ws.Range("C2").Value = 'Input for crop water demands/ET'
$ws.Range("E2").Value = 'Historical climate datasets: Livneh et al (2013); PRISM (2016)'
$ws.Range("A3").Value = 'Demand'
$ws.Range("B3").Value = 'Precipitation'
$ws.Range("C3").Value = 'Input for crop water demands/ET'
$ws.Range("E3").Value = 'Historical climate datasets: Livneh et al (2013); PRISM (2016)'
$ws.Range("A4").Value = 'Demand'
$ws.Range("B4").Value = 'Humidity'
$ws.Range("C4").Value = 'Input for crop water demands/ET'
$ws.Range("E4").Value = 'Historical climate datasets: Livneh et al (2013); PRISM (2016)'
$ws.Range("A5").Value = 'Demand'
$ws.Range("B5").Value = 'Wind speed'
$ws.Range("C5").Value = 'Input for crop water demands/ET'
$ws.Range("E5").Value = 'Historical climate datasets: Livneh et al (2013); PRISM (2016)'
$ws.Range("A6").Value = 'Demand'
$ws.Range("B6").Value = 'Soil water capacity'
$ws.Range("C6").Value = 'Input for crop water demands/ET'
$ws.Range("E6").Value = '?'
$ws.Range("A7").Value = 'Demand'
$ws.Range("B7").Value = 'Soil depth'
$ws.Range("C7").Value = 'Input for crop water demands/ET'
$ws.Range("E7").Value = '?'
$ws.Range("A8").Value = 'Demand'
$ws.Range("B8").Value = 'Planting dates'
$ws.Range("C8").Value = 'Input for crop water demands/ET'
$ws.Range("E8").Value = 'Sacramento - San Joaquin Basin Study (Reclamation, 2014C)'
$ws.Range("A9").Value = 'Demand'
$ws.Range("B9").Value = 'Season length'
$ws.Range("C9").Value = 'Input for crop water demands/ET'
$ws.Range("E9").Value = 'Sacramento - San Joaquin Basin Study (Reclamation, 2014C)'
$ws.Range("A10").Value = 'Demand'
$ws.Range("B10").Value = 'Single crop coefficient'
$ws.Range("C10").Value = 'Input for crop water demands/ET'
$ws.Range("E10").Value = 'Sacramento - San Joaquin Basin Study (Reclamation, 2014C)'
$ws.Range("A11").Value = 'Demand'
$ws.Range("B11").Value = 'Crop specific seasonal application efficiency'
$ws.Range("C11").Value = 'Input to calculate applied water (irrigation water required at the head of field or famr gate)'
$ws.Range("E11").Value = 'Estimated by DWR''s Division of statewide Integrated Water Management'
$ws.Range("A12").Value = 'Demand'
$ws.Range("B12").Value = 'Loss factors'
$ws.Range("C12").Value = 'Input to calculate applied water (irrigation water required at the head of field or farm gate)'
$ws.Range("E12").Value = 'Derived from DWR models and set to 1.0'
$ws.Range("A13").Value = 'Demand'
$ws.Range("B13").Value = 'Area classes'
$ws.Range("C13").Value = 'Input for land use'
$ws.Range("E13").Value = 'California Spatial Information Library (CalSIL); County Land Use Suverys DWR DSIWM (include over 70 crop classifications); County and regional integrated water resources plans and integrated water management plans; Reclamation CVP supply contract renewal and supporting environmental documents'
$ws.Range("A14").Value = 'Demand'
$ws.Range("B14").Value = 'Crop water demands'
$ws.Range("C14").Value = 'daily time step using dual crop coefficient approach described in FAO Drainage Paper (Allen et al. 1998) - MABIA method. Requires inputs of temperature, precipitation, humidity, wind speed. These inputs used to calculated a reference evapotranspiration using the Penman-Monteith Equation. Also requires soil parameters such as soil water capacity and soil depth. The Soil Conservation Service curve number method is used to calculate effective rainfall and rainfall-runoff. '
$ws.Range("D14").Value = 'Daily'
$ws.Range("E14").Value = 'FAO Drainage Paper (Allen et al. 1998) - MABIA method; Crop use parameters based on the Sacramento - San Joaquin Basin Study (Reclamation, 2014C) - planting dates, season length, single crop coefficient'
$ws.Range("B15").Value = 'Climate'
$ws.Range("C15").Value = 'Historical climate data needed for 1922-2015. Two spatially interpolated, gridded datasets developed. One dataset provided daily precipitation, max/min temperature, and wind speed for 1915-2011 on a 1/16 degree grid. the other dataset is combination of daily and monthly data at a 4km grid'
$ws.Range("B16").Value = 'Urban water demands'
$ws.Range("C16").Value = 'historical purveyor production data for 2006-2010 for major cities and towns; based on population data for smaller communities'
$ws.Range("E16").Value = 'DSIWM datasets summarized in California Water Plan (Bulletin 160-09 Series, Bulletin 166 Series), industrial water use reports (Bulletin 124 Series); water use data from 1998-2003 (DWR, 2011). Urban water demans were determined mostly using Public Water System Statistics (PWSS) questionnaires'
$ws.Range("B17").Value = 'South of Delta Demands'
$ws.Range("E17").Value = 'Reclamation''c CVP Contractor data'
$ws.Range("B18").Value = 'Land cover'
$ws.Range("E18").Value = 'National Land Cover Database 2011'

# Row heights
$ws.Rows.Item(1).RowHeight = 17
$ws.Rows.Item(2).RowHeight = 51
$ws.Rows.Item(3).RowHeight = 51
$ws.Rows.Item(4).RowHeight = 51
$ws.Rows.Item(5).RowHeight = 51
$ws.Rows.Item(6).RowHeight = 17
$ws.Rows.Item(7).RowHeight = 17
$ws.Rows.Item(8).RowHeight = 34
$ws.Rows.Item(9).RowHeight = 34
$ws.Rows.Item(10).RowHeight = 34
$ws.Rows.Item(11).RowHeight = 51
$ws.Rows.Item(12).RowHeight = 34
$ws.Rows.Item(13).RowHeight = 204
$ws.Rows.Item(14).RowHeight = 153
$ws.Rows.Item(15).RowHeight = 102
$ws.Rows.Item(16).RowHeight = 187
$ws.Rows.Item(17).RowHeight = 34
$ws.Rows.Item(18).RowHeight = 34

# Selection
$ws.Range("A15").Select()
